$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B width change (match column C's width, which is the target 15.4 char width)
$ws.Columns.Item(2).ColumnWidth = 14.71

# Update cell values in column B
$ws.Range("B2").Value = 866700000.0
$ws.Range("B3").Value = 533300000.0
$ws.Range("B4").Value = -14400000.0
$ws.Range("B5").Value = -432000000.0
$ws.Range("B6").Value = 499300000.0
$ws.Range("B7").Value = -553000000.0
$ws.Range("B8").Value = 41800000.0
$ws.Range("B9").Value = -148900000.0
$ws.Range("B11").Value = 1930200000.0
$ws.Range("B12").Value = -97000000.0
$ws.Range("B13").Value = -643700000.0
$ws.Range("B14").Value = -1400000.0
$ws.Range("B15").Value = -78500000.0
$ws.Range("B16").Value = -739300000.0
$ws.Range("B17").Value = -417200000.0
$ws.Range("B18").Value = 125200000.0
$ws.Range("B19").Value = -357400000.0
$ws.Range("B20").Value = -112100000.0
$ws.Range("B21").Value = -761500000.0
$ws.Range("B22").Value = 298300000.0
$ws.Range("B23").Value = 727700000.0
$ws.Range("B24").Value = 3574300000.0
$ws.Range("B25").Value = 4302000000.0
$ws.Range("B26").Value = 73900000.0
$ws.Range("B27").Value = -357400000.0
$ws.Range("B28").Value = 514899900.0
$ws.Range("B29").Value = -1400000.0
$ws.Range("B30").Value = 125200000.0
